# Daily attendance processing - reorders the "Recorded By" (column G) list
# so that any "System" entries (case-insensitive) appear first, followed by
# the remaining recorder names/emails, preserving their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text.Split(",")
    $systemEntries = @()
    $otherEntries = @()

    foreach ($part in $parts) {
        $trimmed = $part.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemEntries += $trimmed
        } else {
            $otherEntries += $trimmed
        }
    }

    if ($systemEntries.Count -gt 0) {
        $ordered = $systemEntries + $otherEntries
        $newText = $ordered -join ", "
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
